$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
Write-Host $ws.Name
$ws.Range("H11").Value = 1427.091
